$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells so numeric-looking strings
# (e.g. "8.70", "17.21") keep their exact formatting instead of being
# auto-converted to numbers by Excel.
$cells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","E17","D18","D19","E19","D20","E20","D21","E21","D22","E22","E23","E24","D25","E25","E26","D27","E27","D28","E28","D29","E29","E30","D31","E31","E32","E33","E34","D35","E35","E36","E37","D38","E38","E39","B40","C40","D40","E40","B41","C41","D41","E41","D42","E42","E43","D44","E44","D45","E45","D46","E46","E47","E48","E49","D50","E50","D51","E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "35.245.37"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.908.99"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "0.722"
$ws.Range("E5").Value = "  +9.02%  "
$ws.Range("D6").Value = "255.71"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "40.65"
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  +7.26%  "
$ws.Range("D10").Value = "52.86"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "0.0762"
$ws.Range("E11").Value = "  +5.85%  "
$ws.Range("D12").Value = "0.0988"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "2.185.73"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "12.86"
$ws.Range("E14").Value = "  +6.51%  "
$ws.Range("D15").Value = "0.728"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "1.936.86"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "35.260.96"
$ws.Range("D19").Value = "74.87"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "243.87"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "13.05"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("E23").Value = "  +5.44%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +6.86%  "
$ws.Range("E26").Value = "  +4.01%  "
$ws.Range("D27").Value = "166.32"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").Value = "8.70"
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").Value = "18.78"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("D31").Value = "4.129.38"
$ws.Range("E31").Value = "  +19.48%  "
$ws.Range("E32").Value = "  +5.77%  "
$ws.Range("E33").Value = "  +13.73%  "
$ws.Range("E34").Value = "  +22.07%  "
$ws.Range("D35").Value = "0.0589"
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "0.915"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0219"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "17.21"
$ws.Range("E41").Value = "  +5.21%  "
$ws.Range("D42").Value = "96.55"
$ws.Range("E42").Value = "  +7.19%  "
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").Value = "0.0655"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").Value = "1.337.00"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "45.01"
$ws.Range("E50").Value = "  -6.51%  "
$ws.Range("D51").Value = "0.0754"
$ws.Range("E51").Value = "  +6.62%  "
